$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slashes replaced with dashes)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows whose day-of-month is <= 12, where Excel's locale date parser could
# otherwise misread "dd-mm-yyyy" as "mm-dd-yyyy" and silently convert the
# text into a date serial number. Force those specific cells to Text format
# first so the literal dashed string is preserved, just like the others.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($row in ($dates.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($row, 1)
    if ($ambiguousRows -contains $row) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dates[$row]
}

# Value updates for attendance columns (D=4, E=5, G=7, H=8)
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 8).Value = 0

$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 8).Value = 0
